# Add data for 2021-12-22 (workbook tracks "through 12-14" cumulative carjacking arrest stats)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the (only) worksheet and update its title cell / row label
$ws.Name = "Through 2021-12-14"
$ws.Range("A14").Value = "December (through 12-14)"

# Row 14 - "December (through 12-14)" values
$ws.Range("B14").Value = 3
$ws.Range("C14").Value = 12
$ws.Range("D14").Value = 0.2
$ws.Range("E14").Value = 5
$ws.Range("F14").Value = 39
$ws.Range("G14").Value = 0.1136
$ws.Range("H14").Value = 3
$ws.Range("I14").Value = 43
$ws.Range("J14").Value = 0.0652
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 28
$ws.Range("M14").Value = 0.0968
$ws.Range("N14").Value = 3
$ws.Range("O14").Value = 20
$ws.Range("P14").Value = 0.1304
$ws.Range("Q14").Value = 3
$ws.Range("R14").Value = 66
$ws.Range("S14").Value = 0.0435
$ws.Range("T14").Value = 1
$ws.Range("U14").Value = 100
$ws.Range("V14").Value = 0.0099

# New cell V14 needs the same percent format as the other "arrest_rate" columns ("0.0%");
# T14 (arrest_made) keeps the default/general format, matching its column siblings.
$ws.Range("V14").NumberFormat = "0.0%"

# Row 15 - "Total" values
$ws.Range("B15").Value = 36
$ws.Range("C15").Value = 270
$ws.Range("D15").Value = 0.1176
$ws.Range("E15").Value = 65
$ws.Range("F15").Value = 542
$ws.Range("G15").Value = 0.1071
$ws.Range("H15").Value = 66
$ws.Range("I15").Value = 801
$ws.Range("J15").Value = 0.0761
$ws.Range("K15").Value = 77
$ws.Range("L15").Value = 636
$ws.Range("M15").Value = 0.108
$ws.Range("N15").Value = 57
$ws.Range("O15").Value = 500
$ws.Range("P15").Value = 0.1023
$ws.Range("Q15").Value = 67
$ws.Range("R15").Value = 1266
$ws.Range("S15").Value = 0.0503
$ws.Range("T15").Value = 101
$ws.Range("U15").Value = 1643
$ws.Range("V15").Value = 0.0579
